$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7104350185607621
$ws.Range("F2").Value = 0.4741996535665807
$ws.Range("J2").Value = -0.6365504684646589
$ws.Range("B3").Value = 0.8813651726852317
$ws.Range("F3").Value = -0.5107531751786559
$ws.Range("J3").Value = 0.659888126417628
$ws.Range("B4").Value = 0.8284888105101889
$ws.Range("F4").Value = -0.4669612063638096
$ws.Range("J4").Value = 0.6021435811555643
$ws.Range("B5").Value = 0.8166713332729186
$ws.Range("F5").Value = -0.6440512652172118
$ws.Range("J5").Value = -0.9072544033899735
$ws.Range("B6").Value = 0.7856463492376675
$ws.Range("F6").Value = -0.3939408733220401
$ws.Range("J6").Value = -0.6061675427787652
$ws.Range("B7").Value = 0.7621859773132071
$ws.Range("F7").Value = -0.3688413192581249
$ws.Range("J7").Value = -0.592590661221202
$ws.Range("B8").Value = 0.8167642621291933
$ws.Range("F8").Value = 0.5060171112935425
$ws.Range("J8").Value = 0.604040810550924
$ws.Range("B9").Value = 0.7298409743163429
$ws.Range("F9").Value = -0.380952682856165
$ws.Range("J9").Value = -0.545852142808624
$ws.Range("B10").Value = 0.7873838827589105
$ws.Range("F10").Value = 0.3448784394287218
$ws.Range("J10").Value = -0.6201922654759529
$ws.Range("B11").Value = 0.1633157729489279
$ws.Range("F11").Value = -0.4453151504741452
$ws.Range("J11").Value = 0.3231813326375269
$ws.Range("B12").Value = 0.5316860107187826
$ws.Range("F12").Value = 0.3867971054479267
$ws.Range("J12").Value = 0.6149070381258281
$ws.Range("B13").Value = 0.7099187301197926
$ws.Range("F13").Value = -0.2911211630214038
$ws.Range("J13").Value = 0.706155190533714
$ws.Range("B14").Value = 0.6860996232423073
$ws.Range("F14").Value = 0.3970651141398553
$ws.Range("J14").Value = 0.5746657963495071
$ws.Range("B15").Value = 0.7099765534371616
$ws.Range("F15").Value = -0.3072152285082598
$ws.Range("J15").Value = -0.5426382720762762
$ws.Range("B16").Value = 0.2853797359787679
$ws.Range("F16").Value = -0.3796307755022347
$ws.Range("J16").Value = 13.9510825385243
$ws.Range("B17").Value = 0.4272290166072092
$ws.Range("F17").Value = 0.2232305707268276
$ws.Range("J17").Value = -0.5573401192685299
$ws.Range("B18").Value = 0.1348496652117076
$ws.Range("F18").Value = 0.4348570575493335
$ws.Range("J18").Value = -0.6426493068944772
$ws.Range("B19").Value = -0.05145976892022799
$ws.Range("F19").Value = 0.2912410106976744
$ws.Range("J19").Value = -6.997510347546987
$ws.Range("B20").Value = 0.2831366615843744
$ws.Range("F20").Value = 0.2688418847798602
$ws.Range("J20").Value = -0.5494584946267194
$ws.Range("B21").Value = 0.04648977929059236
$ws.Range("F21").Value = -0.5152074001612627
$ws.Range("J21").Value = 0.3778642395909452
